# Rebrand the IT Training Budget Estimates workbook to an
# "Artificial Intelligence and Machine Learning" themed template,
# restoring the original multi-industry template wording.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Instructions & User Guide
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Instructions & User Guide")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning Comprehensive Budget - User Guide & Instructions"
$ws.Range("A56").Value = "📋 ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING PROJECT OVERVIEW"
$ws.Range("B59").Value = "Data Scientists, ML Engineers, AI Architects, DevOps Engineers..."
# Touch E65 (bottom-right corner of the original template's used range)
# with a named style so the sheet's used-range/dimension is restored to
# A1:E65 without introducing a new cell style.
$ws.Range("E65").Style = "Normal"

# ---------------------------------------------------------------
# Budget Summary
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Budget Summary")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Executive Budget Summary"

# ---------------------------------------------------------------
# Resources
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resources")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Resources Budget"
$ws.Range("A4").Value = "Data Scientists"
$ws.Range("A5").Value = "ML Engineers"
$ws.Range("A9").Value = "Business Analysts"

# ---------------------------------------------------------------
# Logistics
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Logistics")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Logistics Budget"

# ---------------------------------------------------------------
# Technology
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Technology")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Technology Budget"

# ---------------------------------------------------------------
# Training
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Training Budget"
$ws.Range("A4").Value = "AI/ML Certification Programs"

# ---------------------------------------------------------------
# Contingency
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contingency")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Contingency Budget"

# ---------------------------------------------------------------
# Timeline
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Timeline")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Budget Timeline"
